$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the Code Review 1 scores (column B) for the four team members
$ws.Range("B3").Value = 25
$ws.Range("B4").Value = 25
$ws.Range("B5").Value = 25
$ws.Range("B6").Value = 25

# Hide the helper/aggregate columns C:F (Code Review 2/3, Final Deliverable, names)
$ws.Range("F1").ColumnWidth = 0
$ws.Range("C:F").EntireColumn.Hidden = $true

# Move the active selection to K13 (matches the reviewer's cursor position)
$ws.Range("K13").Select()
